# #5: property boat&car done
# Fix the "car" (汽車) sheet: row 1 was a stray duplicate of row 2's data
# instead of real column headers, and the row only carried the first six
# columns. Turn row 1 into proper headers, add a "capacity" column, and
# extend both rows out to the full 13-column schema used by the other
# property sheets (name, capacity, owner, register_date, register_reason,
# acquire_value, property_category, category, date, legislator_name,
# legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車 (car)

# --- Row 1: turn the leftover data row into real headers ------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"

# --- Extend row 1 (headers) and row 2 (data) out to column N --------------
# Copy existing formatting (border/bold/alignment for row1, plain for row2)
# into the new cells first, so the new cells keep the sheet's look.
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

# New header cells
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# New data cells (row 2), matching the land/building sheets' constant columns.
# The register/acquire date strings in this workbook are plain text (not real
# dates), so force text with a leading apostrophe to stop Excel from
# reinterpreting "2011-12-26" as a date serial number.
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"
$ws.Cells.Item(2, 10).Value = "'2011-12-26"
$ws.Cells.Item(2, 11).Value = "潘孟安"
$ws.Cells.Item(2, 12).Value = 1376
$ws.Cells.Item(2, 13).Value = "tmp6a821"
$ws.Cells.Item(2, 14).Value = 34

# Re-stamp plain formatting over row 2's new cells: entering "'2011-12-26"
# stuck a text number format on J2, and this keeps every new cell visually
# consistent with the rest of the (unformatted) data row.
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

# Clear the clipboard marquee left behind by Copy()
$excel.CutCopyMode = $false
